$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "id"
$ws.Range("C1").Value = "address"
$ws.Range("F1").Value = "record"

$ws.Range("F1").Select()
